# listaPERFUMARIA.xlsx - "As classes foram devidamente separadas"
# Converts row 2's Meta/Meta.AC/Venda/Venda.AC/Sobras/P columns to real numbers,
# inserts a new numeric row 3, and a new text-formatted row 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: change B2:G2 from text to real numeric values ---
$ws.Cells.Item(2, 2).Value = 1025      # B2 Meta
$ws.Cells.Item(2, 3).Value = 1025      # C2 Meta.AC
$ws.Cells.Item(2, 4).Value = 1040      # D2 Venda
$ws.Cells.Item(2, 5).Value = 1040      # E2 Venda.AC
$ws.Cells.Item(2, 6).Value = 15        # F2 Sobras
$ws.Cells.Item(2, 7).Value = 101.46    # G2 P

# --- Row 3 (new): numeric class row ---
$ws.Cells.Item(3, 1).Value = "31/08/2000"
$ws.Cells.Item(3, 2).Value = 2000
$ws.Cells.Item(3, 3).Value = 3025
$ws.Cells.Item(3, 4).Value = 2000
$ws.Cells.Item(3, 5).Value = 3040
$ws.Cells.Item(3, 6).Value = 15
$ws.Cells.Item(3, 7).Value = 100.5

# --- Row 4 (new): text class row (values kept as literal text, e.g. "2000.00") ---
$ws.Cells.Item(4, 1).Value = "31/02/2000"

$textCells = @(2, 3, 4, 5, 6, 7)
$textValues = @("2000.00", "5025.00", "2000.00", "5040.00", "15.00", "100.30")

for ($i = 0; $i -lt $textCells.Length; $i++) {
    $cell = $ws.Cells.Item(4, $textCells[$i])
    $cell.NumberFormat = "@"
    $cell.Value = $textValues[$i]
    $cell.Style = "Normal"
}
